# ~ arrow letter text size change
#
# Slide 2 has two red "Up Arrow" shapes pointing at a screenshot, labelled
# "A" (drawn inside the first arrow) and "B" (a separate textbox next to the
# second, rotated arrow). Both arrows are enlarged and the two labels grow
# from 14pt to 28pt to match.
#
# NOTE on units: Shape.Left/Top/Width/Height are in points (1 pt = 12700
# EMU) in the PowerPoint object model, not EMU and not inches. The literal
# point values below are chosen so that they reproduce the exact target EMU
# offsets/extents once the host round-trips them through its internal
# (single-precision) storage.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# Shape "Arrow: Up 7" (the upright arrow labelled "A")
#   off  5930495,3189798  -> 5694933,3002555 EMU
#   ext   472535, 953985  ->  802133,1619399 EMU
# ---------------------------------------------------------------------
$arrowA = $s.Shapes.Item(2)
$arrowA.Left   = 448.4199212598425
$arrowA.Top    = 236.4216538433071
$arrowA.Width  = 63.16007874015748
$arrowA.Height = 127.51173228346457
$arrowA.TextFrame.TextRange.Font.Size = 28

# ---------------------------------------------------------------------
# Shape "Arrow: Up 9" -> renamed "Arrow: Up 12" (the rotated arrow; its
# text placeholder is empty, formatting lives only in endParaRPr)
#   off  3485958,5874974  -> 3095667,5399289 EMU
#   ext   472535, 953985  ->  802133,1619399 EMU
# ---------------------------------------------------------------------
$arrowB = $s.Shapes.Item(3)
$arrowB.Left   = 243.75330708661417
$arrowB.Top    = 425.14086614173226
$arrowB.Width  = 63.16007874015748
$arrowB.Height = 127.51173228346457
$arrowB.TextFrame.TextRange.Font.Size = 28
$arrowB.Name = "Arrow: Up 12"

# ---------------------------------------------------------------------
# Shape "TextBox 10" (holds the "B" label next to the rotated arrow)
#   off  3496344,6198077  -> 3206414,5947378 EMU
#   ext   290319, 307777  ->  290319, 523220 EMU (height only)
# ---------------------------------------------------------------------
$textB = $s.Shapes.Item(4)
$textB.Left   = 252.4735489070866
$textB.Top    = 468.29748031496064
$textB.Width  = 22.859763779527558
$textB.Height = 41.198425196850394
$textB.TextFrame.TextRange.Font.Size = 28
